# Apply updated loading_percent values (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 5.243275410925614
$ws.Range("D2").Value = 9.63159552473412
$ws.Range("E2").Value = 9.745651487040259
$ws.Range("F2").Value = 98.45660804162202
$ws.Range("G2").Value = 3.982324374123487
$ws.Range("I2").Value = 73.95804897221333
$ws.Range("J2").Value = 12.81920940909182
$ws.Range("L2").Value = 10.78354274999302
$ws.Range("N2").Value = 18.99769123883691
# Row 3
$ws.Range("C3").Value = 5.104302578158257
$ws.Range("D3").Value = 9.507952370026224
$ws.Range("E3").Value = 9.769517484570418
$ws.Range("F3").Value = 97.91541962454865
$ws.Range("G3").Value = 3.995125036841381
$ws.Range("I3").Value = 73.60855874301286
$ws.Range("J3").Value = 12.86737087250449
$ws.Range("L3").Value = 10.82757635981427
$ws.Range("N3").Value = 18.40031508502705
# Row 4
$ws.Range("C4").Value = 5.019783746689452
$ws.Range("D4").Value = 9.435544475530797
$ws.Range("E4").Value = 9.784905445944011
$ws.Range("F4").Value = 97.61875603217254
$ws.Range("G4").Value = 4.003331750578131
$ws.Range("I4").Value = 73.42041106809661
$ws.Range("J4").Value = 12.89946642057427
$ws.Range("L4").Value = 10.85617237651317
$ws.Range("N4").Value = 18.02485520896364
# Row 5
$ws.Range("C5").Value = 4.985594426278215
$ws.Range("D5").Value = 9.406939693557336
$ws.Range("E5").Value = 9.791361987479331
$ws.Range("F5").Value = 97.50679875096142
$ws.Range("G5").Value = 4.0067642038305
$ws.Range("I5").Value = 73.35035957527099
$ws.Range("J5").Value = 12.91317891532721
$ws.Range("L5").Value = 10.86821925778719
$ws.Range("N5").Value = 17.86990355188116
# Row 6
$ws.Range("C6").Value = 4.979934027558234
$ws.Range("D6").Value = 9.402244896250078
$ws.Range("E6").Value = 9.792445349846481
$ws.Range("F6").Value = 97.48874706387885
$ws.Range("G6").Value = 4.007339509630253
$ws.Range("I6").Value = 73.33912641365839
$ws.Range("J6").Value = 12.91549407390703
$ws.Range("L6").Value = 10.8702434657943
$ws.Range("N6").Value = 17.84406337566596
# Row 7
$ws.Range("C7").Value = 5.019321568820033
$ws.Range("D7").Value = 9.435155025069633
$ws.Range("E7").Value = 9.784991767126741
$ws.Range("F7").Value = 97.61720999822766
$ws.Range("G7").Value = 4.003377683767832
$ws.Range("I7").Value = 73.41943956832101
$ws.Range("J7").Value = 12.89964878968791
$ws.Range("L7").Value = 10.85633324833795
$ws.Range("N7").Value = 18.02277304766463
# Row 8
$ws.Range("C8").Value = 5.195218212051261
$ws.Range("D8").Value = 9.588247358543605
$ws.Range("E8").Value = 9.753728849786659
$ws.Range("F8").Value = 98.2626014399999
$ws.Range("G8").Value = 3.986666547945268
$ws.Range("I8").Value = 73.83204242442777
$ws.Range("J8").Value = 12.83529069276941
$ws.Range("L8").Value = 10.79840295209613
$ws.Range("N8").Value = 18.79364780656866
# Row 9
$ws.Range("C9").Value = 5.544512110489636
$ws.Range("D9").Value = 9.915351719947468
$ws.Range("E9").Value = 9.698191484204195
$ws.Range("F9").Value = 99.81138231629886
$ws.Range("G9").Value = 3.95660891197088
$ws.Range("I9").Value = 74.85184375658956
$ws.Range("J9").Value = 12.72918208221869
$ws.Range("L9").Value = 10.69709614158877
$ws.Range("N9").Value = 20.2273683202997
# Row 10
$ws.Range("C10").Value = 5.801171664348588
$ws.Range("D10").Value = 10.17068383713315
$ws.Range("E10").Value = 9.660825974703046
$ws.Range("F10").Value = 101.1218036377929
$ws.Range("G10").Value = 3.936120989532287
$ws.Range("I10").Value = 75.73024894882062
$ws.Range("J10").Value = 12.66357806392304
$ws.Range("L10").Value = 10.63005304496841
$ws.Range("N10").Value = 21.22223697909767
# Row 11
$ws.Range("C11").Value = 5.917438816810114
$ws.Range("D11").Value = 10.28978353325449
$ws.Range("E11").Value = 9.644556442143518
$ws.Range("F11").Value = 101.7551583577385
$ws.Range("G11").Value = 3.927133830118985
$ws.Range("I11").Value = 76.15789633293703
$ws.Range("J11").Value = 12.63644007479575
$ws.Range("L11").Value = 10.60113406152992
$ws.Range("N11").Value = 21.66018057919901
# Row 12
$ws.Range("C12").Value = 5.96135840086671
$ws.Range("D12").Value = 10.33527996911248
$ws.Range("E12").Value = 9.638498886602722
$ws.Range("F12").Value = 102.0003087632882
$ws.Range("G12").Value = 3.923777404985022
$ws.Range("I12").Value = 76.32385627121387
$ws.Range("J12").Value = 12.6265549217905
$ws.Range("L12").Value = 10.59040842936012
$ws.Range("N12").Value = 21.82377585681918
# Row 13
$ws.Range("C13").Value = 5.951904915395414
$ws.Range("D13").Value = 10.32546433052
$ws.Range("E13").Value = 9.639798914351918
$ws.Range("F13").Value = 101.9472757098005
$ws.Range("G13").Value = 3.924498206141211
$ws.Range("I13").Value = 76.28793538448896
$ws.Range("J13").Value = 12.62866641910782
$ws.Range("L13").Value = 10.59270839106621
$ws.Range("N13").Value = 21.78864458690801
# Row 14
$ws.Range("C14").Value = 5.92105448239006
$ws.Range("D14").Value = 10.29351877632695
$ws.Range("E14").Value = 9.644056019713364
$ws.Range("F14").Value = 101.7752208246568
$ws.Range("G14").Value = 3.926856762878583
$ws.Range("I14").Value = 76.17146951648186
$ws.Range("J14").Value = 12.63561895879289
$ws.Range("L14").Value = 10.60024715078318
$ws.Range("N14").Value = 21.67368539489599
# Row 15
$ws.Range("C15").Value = 5.902142534163306
$ws.Range("D15").Value = 10.27400194071092
$ws.Range("E15").Value = 9.646677039917954
$ws.Range("F15").Value = 101.6705227220417
$ws.Range("G15").Value = 3.928307512086111
$ws.Range("I15").Value = 76.10065366282461
$ws.Range("J15").Value = 12.63992863608554
$ws.Range("L15").Value = 10.60489415541532
$ws.Range("N15").Value = 21.60297336126124
# Row 16
$ws.Range("C16").Value = 5.79356017266661
$ws.Range("D16").Value = 10.16295732205425
$ws.Range("E16").Value = 9.661903765552371
$ws.Range("F16").Value = 101.0811583391496
$ws.Range("G16").Value = 3.936714931062028
$ws.Range("I16").Value = 75.70286458618048
$ws.Range("J16").Value = 12.66540623682011
$ws.Range("L16").Value = 10.63197459920526
$ws.Range("N16").Value = 21.19330956972086
# Row 17
$ws.Range("C17").Value = 5.726795134269771
$ws.Range("D17").Value = 10.09557019533209
$ws.Range("E17").Value = 9.671430429916191
$ws.Range("F17").Value = 100.7291127173018
$ws.Range("G17").Value = 3.941957135118653
$ws.Range("I17").Value = 75.46601090799805
$ws.Range("J17").Value = 12.68173051875504
$ws.Range("L17").Value = 10.64899082976388
$ws.Range("N17").Value = 20.93814219015166
# Row 18
$ws.Range("C18").Value = 5.688349722356528
$ws.Range("D18").Value = 10.0570900526706
$ws.Range("E18").Value = 9.676978542036148
$ws.Range("F18").Value = 100.5301334730141
$ws.Range("G18").Value = 3.945003711349473
$ws.Range("I18").Value = 75.33242109138828
$ws.Range("J18").Value = 12.69137432942896
$ws.Range("L18").Value = 10.65892686745585
$ws.Range("N18").Value = 20.79000725568141
# Row 19
$ws.Range("C19").Value = 5.675326436906881
$ws.Range("D19").Value = 10.04411013268202
$ws.Range("E19").Value = 9.678868864969369
$ws.Range("F19").Value = 100.4633661729463
$ws.Range("G19").Value = 3.94604065407364
$ws.Range("I19").Value = 75.28764403130893
$ws.Range("J19").Value = 12.69468320196416
$ws.Range("L19").Value = 10.66231664001299
$ws.Range("N19").Value = 20.73962067985785
# Row 20
$ws.Range("C20").Value = 5.733907219366929
$ws.Range("D20").Value = 10.10271498144336
$ws.Range("E20").Value = 9.670409207343022
$ws.Range("F20").Value = 100.7662257780586
$ws.Range("G20").Value = 3.941395851380396
$ws.Range("I20").Value = 75.49095102481314
$ws.Range("J20").Value = 12.67996641312099
$ws.Range("L20").Value = 10.64716404023999
$ws.Range("N20").Value = 20.96544799483449
# Row 21
$ws.Range("C21").Value = 5.930119226101021
$ws.Range("D21").Value = 10.30289143570847
$ws.Range("E21").Value = 9.64280281036833
$ws.Range("F21").Value = 101.8256137056617
$ws.Range("G21").Value = 3.926162735560506
$ws.Range("I21").Value = 76.20556943867992
$ws.Range("J21").Value = 12.63356618702003
$ws.Range("L21").Value = 10.59802673204518
$ws.Range("N21").Value = 21.70751365554064
# Row 22
$ws.Range("C22").Value = 6.057708953831852
$ws.Range("D22").Value = 10.43601350788469
$ws.Range("E22").Value = 9.625362264286322
$ws.Range("F22").Value = 102.5489264674501
$ws.Range("G22").Value = 3.916479385065543
$ws.Range("I22").Value = 76.69602939196372
$ws.Range("J22").Value = 12.60552361654176
$ws.Range("L22").Value = 10.56722540899465
$ws.Range("N22").Value = 22.17935961385674
# Row 23
$ws.Range("C23").Value = 5.989682680815246
$ws.Range("D23").Value = 10.36476292712577
$ws.Range("E23").Value = 9.634615999907702
$ws.Range("F23").Value = 102.1600656514109
$ws.Range("G23").Value = 3.921623006817347
$ws.Range("I23").Value = 76.43212575397519
$ws.Range("J23").Value = 12.62028078560847
$ws.Range("L23").Value = 10.5835451161513
$ws.Range("N23").Value = 21.92877110911181
# Row 24
$ws.Range("C24").Value = 5.730692033712003
$ws.Range("D24").Value = 10.09948400665167
$ws.Range("E24").Value = 9.67087068048717
$ws.Range("F24").Value = 100.7494363186425
$ws.Range("G24").Value = 3.941649505628837
$ws.Range("I24").Value = 75.47966755916451
$ws.Range("J24").Value = 12.68076315946363
$ws.Range("L24").Value = 10.64798945453067
$ws.Range("N24").Value = 20.95310750188673
# Row 25
$ws.Range("C25").Value = 5.449832611648653
$ws.Range("D25").Value = 9.824120275499064
$ws.Range("E25").Value = 9.712606294651851
$ws.Range("F25").Value = 99.36204316488789
$ws.Range("G25").Value = 3.964455808103859
$ws.Range("I25").Value = 74.55329849326844
$ws.Range("J25").Value = 12.75572558983968
$ws.Range("L25").Value = 10.72319724458162
$ws.Range("N25").Value = 19.84905939529497
